$d = $word.ActiveDocument

# Pass 1: replace each original value with a unique placeholder token
# to avoid cross-cell substring collisions between old/new values.
# Pass 2: replace each placeholder with its final value.
$pairs = @(
    ,@("2024-07-13 Saturday", "2024-07-14 Sunday", "@@0@@")
    ,@("68-42=", "86-35=", "@@1@@")
    ,@("16+56=", "45-18=", "@@2@@")
    ,@("14+28=", "48-21=", "@@3@@")
    ,@("72-28=", "21+42=", "@@4@@")
    ,@("79-48=", "84-67=", "@@5@@")
    ,@("81-11=", "19+78=", "@@6@@")
    ,@("82-27=", "47-45=", "@@7@@")
    ,@("48+37=", "59-50=", "@@8@@")
    ,@("94-20=", "33+23=", "@@9@@")
    ,@("85-66=", "35+55=", "@@10@@")
    ,@("22+14=", "46+17=", "@@11@@")
    ,@("99-27=", "75-44=", "@@12@@")
    ,@("92-10=", "77-67=", "@@13@@")
    ,@("67+4=", "87-76=", "@@14@@")
    ,@("24-11=", "69-14=", "@@15@@")
    ,@("42+38=", "51-15=", "@@16@@")
    ,@("36+55=", "33+52=", "@@17@@")
    ,@("64-11=", "74-49=", "@@18@@")
    ,@("56-34=", "0+80=", "@@19@@")
    ,@("84-36=", "74-72=", "@@20@@")
    ,@("41+16=", "66+14=", "@@21@@")
    ,@("40+12=", "73-37=", "@@22@@")
    ,@("78-51=", "19+45=", "@@23@@")
    ,@("50+18=", "54-1=", "@@24@@")
    ,@("91-4=", "43+52=", "@@25@@")
    ,@("86-57=", "11+11=", "@@26@@")
    ,@("74-7=", "50-13=", "@@27@@")
    ,@("66-45=", "86-31=", "@@28@@")
    ,@("37-7=", "42+28=", "@@29@@")
    ,@("87-11=", "61+10=", "@@30@@")
    ,@("38+33=", "42+17=", "@@31@@")
    ,@("7-5=", "82+9=", "@@32@@")
    ,@("95-25=", "79-45=", "@@33@@")
    ,@("51+8=", "51-13=", "@@34@@")
    ,@("10+74=", "31+16=", "@@35@@")
    ,@("86-84=", "77-19=", "@@36@@")
    ,@("27+69=", "72-70=", "@@37@@")
    ,@("58-38=", "72+3=", "@@38@@")
    ,@("8+91=", "65+12=", "@@39@@")
    ,@("66-15=", "55-4=", "@@40@@")
    ,@("10+39=", "80-8=", "@@41@@")
    ,@("94-22=", "56-30=", "@@42@@")
    ,@("47-9=", "25+51=", "@@43@@")
    ,@("20-19=", "35-19=", "@@44@@")
    ,@("34+34=", "78-60=", "@@45@@")
    ,@("60-4=", "84-80=", "@@46@@")
    ,@("55+25=", "22+66=", "@@47@@")
    ,@("3+3=", "53-34=", "@@48@@")
    ,@("98-11=", "54-42=", "@@49@@")
    ,@("2+28=", "59-42=", "@@50@@")
    ,@("22+76=", "45+3=", "@@51@@")
    ,@("14+67=", "27+55=", "@@52@@")
    ,@("29-20=", "63-12=", "@@53@@")
    ,@("85-56=", "90-78=", "@@54@@")
    ,@("15+25=", "53-26=", "@@55@@")
    ,@("14-2=", "68+26=", "@@56@@")
    ,@("91-2=", "71-24=", "@@57@@")
    ,@("48-39=", "31+65=", "@@58@@")
    ,@("42-29=", "76+5=", "@@59@@")
    ,@("61+34=", "64-49=", "@@60@@")
    ,@("54-12=", "81-8=", "@@61@@")
    ,@("42+5=", "35+63=", "@@62@@")
    ,@("88-15=", "1-1=", "@@63@@")
    ,@("37+51=", "0+76=", "@@64@@")
    ,@("98-43=", "94-57=", "@@65@@")
    ,@("66-27=", "91-83=", "@@66@@")
    ,@("77-53=", "53-20=", "@@67@@")
    ,@("91-23=", "18+9=", "@@68@@")
    ,@("78+20=", "70-50=", "@@69@@")
    ,@("43+36=", "59-9=", "@@70@@")
    ,@("88-10=", "23-14=", "@@71@@")
    ,@("49-44=", "89-72=", "@@72@@")
    ,@("88-41=", "74-4=", "@@73@@")
    ,@("59+7=", "12+27=", "@@74@@")
    ,@("40+52=", "25+70=", "@@75@@")
    ,@("48-16=", "79-1=", "@@76@@")
    ,@("76-44=", "88-82=", "@@77@@")
    ,@("45-17=", "4+34=", "@@78@@")
    ,@("0+41=", "5+78=", "@@79@@")
    ,@("84-65=", "99-58=", "@@80@@")
    ,@("80-26=", "38+20=", "@@81@@")
    ,@("5+54=", "2+29=", "@@82@@")
    ,@("2+27=", "52-52=", "@@83@@")
    ,@("75-24=", "28+57=", "@@84@@")
    ,@("64-24=", "45+37=", "@@85@@")
    ,@("87-75=", "47+43=", "@@86@@")
    ,@("93-47=", "13+67=", "@@87@@")
    ,@("25+7=", "58-11=", "@@88@@")
    ,@("74+24=", "60-5=", "@@89@@")
    ,@("17+28=", "38-20=", "@@90@@")
    ,@("30-20=", "64-1=", "@@91@@")
    ,@("96-25=", "6+8=", "@@92@@")
    ,@("43-22=", "34+31=", "@@93@@")
    ,@("75-21=", "91-26=", "@@94@@")
    ,@("97+0=", "33+2=", "@@95@@")
    ,@("23+22=", "94-41=", "@@96@@")
    ,@("13+17=", "40+10=", "@@97@@")
    ,@("37+35=", "15+18=", "@@98@@")
    ,@("21+69=", "66-63=", "@@99@@")
    ,@("38+21=", "18+16=", "@@100@@")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $token = $pair[2]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $token, 2) | Out-Null
}

foreach ($pair in $pairs) {
    $token = $pair[2]
    $new = $pair[1]
    $d.Content.Find.Execute($token, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
